$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 3 through 34,
# replacing the previous "Strike#"-derived values.
$newValues = @{
    3  = 1
    4  = 2
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 2
    10 = 1
    11 = 2
    12 = 3
    13 = 2
    14 = 5
    15 = 3
    16 = 2
    17 = 3
    18 = 2
    19 = 0
    20 = 3
    21 = 3
    22 = 1
    23 = 0
    24 = 3
    25 = 0
    26 = 2
    27 = 0
    28 = 3
    29 = 1
    30 = 2
    31 = 2
    32 = 3
    33 = 3
    34 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
